# moved registration to seperate page
# Slide 25: "Rounded Rectangle 10" code block - recolor part of the url(...) line
# to a flat gray (srgbClr BFBFBF) instead of the themed bg1/lumMod(95%), and
# merge the trailing "," into the preceding run (removing the separate comma run).
# Also Slide 25: "Rounded Rectangle 13" - drop the stray "r" before the smart
# quote in " <- matches r'^$'" -> " <- matches '^$'".

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(25)

$BFBFBF = 12566463  # 0xBFBFBF as a VBA-style RGB long (R + G*256 + B*65536)

# ---- Shape "Rounded Rectangle 10": url(...) line ----
$code = $s.Shapes.Item(6)
$tr = $code.TextFrame.TextRange

# Remove the separate trailing "," run (originally chars 128) first, working
# right-to-left so earlier character offsets stay valid.
$commaRun = $tr.Characters(128, 1)
$commaRun.Text = ""

# Fold the comma into the previous run's text and recolor it (recolor first,
# then grow the text - setting Text last keeps the edit applied).
$homeRun = $tr.Characters(113, 15)
$homeRun.Font.Color.RGB = $BFBFBF
$homeRun.Text = "', name='home'),"

$eventRun = $tr.Characters(97, 16)
$eventRun.Font.Color.RGB = $BFBFBF

$parenRun = $tr.Characters(88, 9)
$parenRun.Font.Color.RGB = $BFBFBF

$urlRun = $tr.Characters(85, 3)
$urlRun.Font.Color.RGB = $BFBFBF

# ---- Shape "Rounded Rectangle 13": matches caption ----
# Drop the stray "r" before the curly single quotes (U+2019), i.e.
# " <- matches r'^$'" -> " <- matches '^$'" (quotes are typographic, not ASCII).
$caption = $s.Shapes.Item(9)
$ctr = $caption.TextFrame.TextRange
$matchRun = $ctr.Characters(25, 17)
$matchRun.Text = " <- matches " + [char]0x2019 + "^`$" + [char]0x2019
